# "Fruta / hortaliza, semanal" weekly refresh.
#
# A new weekly observation is inserted at row 13 (date 2022-09-26 /
# serial 44830), which pushes the existing rows 13-20 down by one
# (they become rows 14-21), growing the used range from A1:R20 to
# A1:R21.
#
# Implemented as: copy whole rows bottom-up (20->21, 19->20, ... 13->14)
# so nothing is clobbered before it is read, then overwrite row 13 with
# the brand-new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 20; $r -ge 13; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}
$excel.CutCopyMode = $false

# New record for row 13 (Volumen/Precio columns change; everything else
# in the row - market, region, category, unit, origin, etc. - stays the
# same as its neighbours, so only these cells need to be written).
$ws.Range("D13").Value = 44830
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 12000
$ws.Range("P13").Value = 800
